# Auto update Excel log
# Appends newly-logged sensor rows to the ALERTS, PIR, Humidity and
# Temperature sheets (the "2026-01-28 15:35/15:36" batch of readings that
# was captured after the previous last row on each sheet).

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param([string]$SheetName, [int]$StartRow, [object[]]$Rows, [string[]]$TextColumns)

    $ws = $wb.Worksheets.Item($SheetName)

    $endRow = $StartRow + $Rows.Count - 1

    # Column A always holds "yyyy-mm-dd" style text in this log, and some
    # sheets (e.g. Humidity's "NN.N%" readings) hold other text that looks
    # numeric. Force Text format on those columns first so Excel doesn't
    # auto-convert the literal strings into date serials / percentages.
    foreach ($col in $TextColumns) {
        $rangeStr = $col + $StartRow + ":" + $col + $endRow
        $ws.Range($rangeStr).NumberFormat = "@"
    }

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $rowValues = $Rows[$i]
        for ($c = 0; $c -lt $rowValues.Count; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
        }
    }

    # The temporary Text format was only needed to stop Excel reinterpreting
    # the literal strings above as dates/percentages; drop it again so the
    # new rows pick up the sheet's normal (unstyled) look, same as every
    # other data row already on the sheet.
    if ($TextColumns.Count -gt 0) {
        $rangeStr = "A" + $StartRow + ":F" + $endRow
        $ws.Range($rangeStr).ClearFormats()
    }
}

# ---------------------------------------------------------------------
# ALERTS — 3 new threshold-breach alerts
# ---------------------------------------------------------------------
$alertsRows = @(
    @("2026-01-28", "15:35:13", "15:00", "Bathroom", "MINIMAL",  "MINIMAL ALERT: Bathroom occupied, no motion > 20s."),
    @("2026-01-28", "15:35:33", "15:00", "Bathroom", "MODERATE", "MODERATE ALERT: Bathroom occupied, no motion > 40s."),
    @("2026-01-28", "15:35:54", "15:00", "Bathroom", "CRITICAL", "CRITICAL ALERT: Bathroom occupied, no motion > 60s.")
)
Add-LogRows "ALERTS" 2 $alertsRows @("A")

# ---------------------------------------------------------------------
# PIR — 13 new "No Motion" / Inactive readings
# ---------------------------------------------------------------------
$pirTimes = @(
    "15:35:03","15:35:05","15:35:10","15:35:15","15:35:20","15:35:25","15:35:30",
    "15:35:35","15:35:40","15:35:45","15:35:50","15:35:55","15:36:00"
)
$pirRows = @()
foreach ($t in $pirTimes) {
    $pirRows += ,@("2026-01-28", $t, "15:00", "Bathroom", "No Motion", "Inactive")
}
Add-LogRows "PIR" 81 $pirRows @("A")

# ---------------------------------------------------------------------
# Humidity — 12 new readings (Active)
# ---------------------------------------------------------------------
$humidityData = @(
    @("15:35:02","88.6%"),
    @("15:35:06","88.5%"),
    @("15:35:10","88.4%"),
    @("15:35:18","88.4%"),
    @("15:35:22","88.4%"),
    @("15:35:27","87.5%"),
    @("15:35:31","88.4%"),
    @("15:35:35","87.5%"),
    @("15:35:43","88.4%"),
    @("15:35:47","87.5%"),
    @("15:35:55","87.5%"),
    @("15:35:59","88.4%")
)
$humidityRows = @()
foreach ($pair in $humidityData) {
    $humidityRows += ,@("2026-01-28", $pair[0], "15:00", "Bathroom", $pair[1], "Active")
}
Add-LogRows "Humidity" 83 $humidityRows @("A", "E")

# ---------------------------------------------------------------------
# Temperature — 12 new readings (Active)
# ---------------------------------------------------------------------
$temperatureData = @(
    @("15:35:02","22.9C"),
    @("15:35:07","22.9C"),
    @("15:35:11","22.9C"),
    @("15:35:19","22.9C"),
    @("15:35:23","22.9C"),
    @("15:35:27","22.9C"),
    @("15:35:31","22.9C"),
    @("15:35:35","22.9C"),
    @("15:35:43","22.9C"),
    @("15:35:47","22.9C"),
    @("15:35:55","22.9C"),
    @("15:35:59","22.9C")
)
$temperatureRows = @()
foreach ($pair in $temperatureData) {
    $temperatureRows += ,@("2026-01-28", $pair[0], "15:00", "Bathroom", $pair[1], "Active")
}
Add-LogRows "Temperature" 83 $temperatureRows @("A")
